# "Generate Report for Handoff"
# The b.md row ("Ready for handoff" row) on each sheet gets refreshed with
# new status / handoff-file / handoff-datetime values reflecting a new
# handoff commit (63290e5768f688058c7b37413b0a5c26c308f864) for b.md.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row for b.md (row 3) gets the new status + datetime.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-23 20:37:38"

# ---------------------------------------------------------------------
# zh-cn sheet: row for b.md (row 3) - Status / Latest Handoff File /
# Latest Handoff Datetime change; the hyperlink display text for the
# Latest Handoff File cell (D3) must point at the new handoff filename.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-23 20:37:34"

# Hyperlinks loaded from the source file cannot be edited/removed one at a
# time through this object model (Delete() on a pre-existing Hyperlink is a
# no-op), so rebuild the sheet's whole Hyperlinks collection: clear it, then
# re-add every entry with its original target + display text, swapping in
# the new display text only for the one that changed (D3).
$zhcn.Range("A1").Hyperlinks.Delete()

$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/10fa69b84d58bd4c57d30a064ff06d03f333fc33/e2e/a.md", "", "", "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1db335bbb86d9e224fba13bf2d5796e616b4bdf6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/9a1d942327e5ddd14a18bc8a9890f10e965c4db6/e2e/a.md", "", "", "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/53c8050c5a6a0b1e86392e8979ffcb23e062a4a5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/10fa69b84d58bd4c57d30a064ff06d03f333fc33/e2e/b.md", "", "", "b.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1db335bbb86d9e224fba13bf2d5796e616b4bdf6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/9a1d942327e5ddd14a18bc8a9890f10e965c4db6/e2e/a.md", "", "", "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/53c8050c5a6a0b1e86392e8979ffcb23e062a4a5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn, but the source file name
# keeps the "a.6631f68b...de-de.xlf" filename for row 2 (a.md), and the
# Latest Handback columns already matched what is new for row 3 (b.md)
# except for file/datetime which change in the same way as zh-cn.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-23 20:37:38"

$dede.Range("A1").Hyperlinks.Delete()

$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/10fa69b84d58bd4c57d30a064ff06d03f333fc33/e2e/a.md", "", "", "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b6a16c9237d1bbfbbcab9b7a685f4d28cc8f765/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/853d15faaa2a86ee69adf5645c4b952dc180faae/e2e/a.md", "", "", "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7c3af6dbdaf883837be921adae3ba2da14b0e5a3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/10fa69b84d58bd4c57d30a064ff06d03f333fc33/e2e/b.md", "", "", "b.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b6a16c9237d1bbfbbcab9b7a685f4d28cc8f765/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/853d15faaa2a86ee69adf5645c4b952dc180faae/e2e/a.md", "", "", "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7c3af6dbdaf883837be921adae3ba2da14b0e5a3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
